$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append new venue/county rows after the existing data (rows 2-99),
# continuing at row 100 through row 108.
$newRows = @(
    @('Bord Na Mona O Connor Park,', 'Offaly'),
    @('Ballina Stephenites', 'Mayo'),
    @('University of Limerick 3G Pitch', 'Limerick'),
    @('Kiltoom', 'Roscommon'),
    @('Eire Og, Carrickmore', 'Tyrone'),
    @('Darver', 'Louth'),
    @('Páirc na hÓige, Maguiresbridge', 'Fermanagh'),
    @("O'Tooles, Dublin", 'Dublin'),
    @('Bellefield, Enniscorthy', 'Wexford')
)

$startRow = 100
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $newRows[$i][0]
    $ws.Cells.Item($r, 2).Value = $newRows[$i][1]
}

# Reset the view's active cell back to the top of the sheet (the saved
# workbook no longer parks the selection on the now-populated A101).
[void]$ws.Range("A1").Select()
